# Update Excel parameters 2
# - Input!C6 (ELKW count): 0 -> 3
# - Input!C7 (WELKW count): 0 -> 3
# - Input!C8 (EPKW flag):   "False" -> "True" (kept as text, not boolean)
# - Row-height refresh on WELKW (rows 3 & 6) and EPKW (row 6), matching the
#   wrap-text reflow Excel performs after the above edits
# - Active sheet/selection moves from EPKW to Input (with EPKW's own
#   selection parked at B14)

$wb = $excel.ActiveWorkbook

$inputSheet = $wb.Worksheets.Item("Input")
$welkw = $wb.Worksheets.Item("WELKW")
$epkw = $wb.Worksheets.Item("EPKW")

# --- Update the Input parameters ---------------------------------------
$inputSheet.Range("C6").Value = 3
$inputSheet.Range("C7").Value = 3

# Write "True" as literal text (not boolean) by computing it via a formula
# and then pasting the result back as a value - this preserves the existing
# cell style and keeps the shared-string text type used throughout the
# sheet instead of Excel's auto bool-coercion.
$inputSheet.Range("C8").Formula = "=""True"""
$inputSheet.Range("C8").Copy()
$inputSheet.Range("C8").PasteSpecial(-4163)

# --- Row-height reflow on dependent sheets ------------------------------
$welkw.Rows.Item(3).RowHeight = 30
$welkw.Rows.Item(6).RowHeight = 45
$epkw.Rows.Item(6).RowHeight = 30

# --- Selection / active sheet changes -----------------------------------
# Select EPKW's new cell first (this sheet loses focus afterwards).
$epkw.Range("B14").Select()

# Finally activate Input and select B3 so it ends up the active tab.
$inputSheet.Activate()
$inputSheet.Range("B3").Select()
